$wb = $excel.ActiveWorkbook

$oldGuid = "ba056d89-61db-4787-85ec-ff51c52bd823"
$newGuid = "64109069-bb76-4707-9289-157a65c12f1c"
$newZhXlf = "$newGuid.fc7598458889aef83d9396afe69e005b0b1cf15e.zh-cn.xlf"
$newDeXlf = "$newGuid.fc7598458889aef83d9396afe69e005b0b1cf15e.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Capture the existing hyperlink target address before touching anything.
$ovAddr = $null
foreach ($hl in $wsOverview.Hyperlinks) {
    $ovAddr = $hl.Address
    break
}

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-20 05:01:58"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $ovAddr, "", "", "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhAddrA = $null
foreach ($hl in $wsZh.Hyperlinks) {
    $zhAddrA = $hl.Address
    break
}

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = "2016-08-20 05:01:53"

# Latest Target File / Latest Handback File become blank (losing their
# hyperlink + hyperlink style in the process).
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""

$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# The sheet only keeps the A2 hyperlink after the edit (I2's is gone).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhAddrA, "", "", "$newGuid.md") | Out-Null

$wsZh.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsZh.Columns.Item(10).ColumnWidth = 20.872143700009268

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deAddrA = $null
foreach ($hl in $wsDe.Hyperlinks) {
    $deAddrA = $hl.Address
    break
}

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = $newDeXlf

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""

$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deAddrA, "", "", "$newGuid.md") | Out-Null

$wsDe.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsDe.Columns.Item(10).ColumnWidth = 20.872143700009268
